# Ultimo commit - prueba tecnica
# Update the "Delete a Member" row (row 7) to target member id 40 instead of 5,
# and move the active selection to H8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("infoAPI")

# EndPoint column (C) and expected response column (H) for the "Delete a Member" row.
$ws.Range("C7").Value = "/api/members/40"
$ws.Range("H7").Value = "Member with id 40 is deleted successfully"

# Update the active cell / selection shown in the sheet view.
$ws.Activate()
$ws.Range("H8").Select()
